$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, pushing existing rows 4-27 down to 5-28.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with this week's data (same template values as the
# rest of the data set, new date + volume/price figures).
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44761
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 100112035
$ws.Cells.Item(4, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 17000
$ws.Cells.Item(4, 12).Value = 18000
$ws.Cells.Item(4, 13).Value = 17400
$ws.Cells.Item(4, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(4, 16).Value = 1160
$ws.Cells.Item(4, 17).Value = 15
$ws.Cells.Item(4, 18).Value = "Hortaliza"
